$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wednesday (row 5) — fill in today's hours worked, matching the style of the
# other fully-populated day rows: a 15-minute-per-column grid of key letters
# (s = sleep, b = break, w = working).
$ws.Range("C5:F5").Value = "s"
$ws.Range("G5:J5").Value = "b"
$ws.Range("K5:P5").Value = "w"
$ws.Range("Q5:R5").Value = "b"
$ws.Range("S5:AL5").Value = "w"

# Recalculate so the COUNTIF/SUM summary cells (rows 13-19) pick up today's
# entries.
$excel.Calculate()

# Leave the selection where the user finished typing.
$ws.Range("K5").Select()
